$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-02 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("35÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=", 2) | Out-Null
$d.Content.Find.Execute("23÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷4=", 2) | Out-Null
$d.Content.Find.Execute("67÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷8=", 2) | Out-Null
$d.Content.Find.Execute("16÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=", 2) | Out-Null
$d.Content.Find.Execute("50÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=", 2) | Out-Null
$d.Content.Find.Execute("32÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("68÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷6=", 2) | Out-Null
$d.Content.Find.Execute("75÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷4=", 2) | Out-Null
$d.Content.Find.Execute("10÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷3=", 2) | Out-Null
$d.Content.Find.Execute("44÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷8=", 2) | Out-Null
$d.Content.Find.Execute("75÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 2) | Out-Null
$d.Content.Find.Execute("53÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷6=", 2) | Out-Null
$d.Content.Find.Execute("92÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=", 2) | Out-Null
$d.Content.Find.Execute("91÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷5=", 2) | Out-Null
$d.Content.Find.Execute("71÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=", 2) | Out-Null
$d.Content.Find.Execute("56÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷5=", 2) | Out-Null
$d.Content.Find.Execute("72÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=", 2) | Out-Null
$d.Content.Find.Execute("81÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=", 2) | Out-Null
$d.Content.Find.Execute("15÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=", 2) | Out-Null
$d.Content.Find.Execute("33÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷9=", 2) | Out-Null
$d.Content.Find.Execute("44÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 2) | Out-Null
$d.Content.Find.Execute("23÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷2=", 2) | Out-Null
$d.Content.Find.Execute("32÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷4=", 2) | Out-Null
$d.Content.Find.Execute("55÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷4=", 2) | Out-Null
$d.Content.Find.Execute("90÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷5=", 2) | Out-Null
